$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: I1 = "I0", J1 = "IF" ---
# Copy formatting from the existing header cell (H1) so the new header
# cells pick up the same bold / bordered / centered style used by the
# other header cells, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Copy()
$ws.Range("J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows 2-67: new columns I (I0) and J (IF) ---
$colI = @(7,6,7,7,5,1,6,4,7,7,6,8,4,1,6,4,5,6,7,7,8,7,6,8,9,7,5,7,6,3,7,7,8,8,6,6,9,8,8,8,7,8,5,5,9,6,8,4,10,7,5,8,7,7,7,9,7,6,5,8,8,4,5,9,6,4)
$colJ = @(7,6,7,8,5,1,7,5,7,7,6,8,4,2,6,5,5,6,7,7,8,7,6,8,9,7,5,7,6,3,7,7,8,8,6,6,9,8,8,8,7,8,5,5,9,6,8,4,10,7,5,8,7,7,7,9,7,6,5,8,8,4,5,9,6,4)

for ($i = 0; $i -lt $colI.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value  = $colI[$i]
    $ws.Cells.Item($row, 10).Value = $colJ[$i]
}

Write-Host "Added columns I (I0) and J (IF) for rows 1-67"
